$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "Files" tab (row 4) Cypher query is rewritten to a trimmed-down
# version of itself that drops the "File Type" and "Breed" output columns
# (object repository revisited: Files tab no longer returns those fields).
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE f.file_format IN ["rtf"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# The active selection moved from C4 to B4 (the cell that was just edited).
$ws.Range("B4").Select()
